$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Authors column (E) for rows 2-4 had their shared-string values replaced
# with cleaned-up copies (one additional level of indentation/whitespace was
# introduced between the comma-separated author entries, consistent with the
# progressive pattern already present in the workbook's shared strings table).

$ws.Range("E2").Value = "[Luciano%Gattinoni%gattinoniluciano@gmail.com%2,               Davide%Chiumello%NULL%3,               Sandra%Rossi%NULL%3]"
$ws.Range("E3").Value = "[Luciano%Gattinoni%NULL%0,               Silvia%Coppola%NULL%2,               Silvia%Coppola%NULL%0,               Massimo%Cressoni%NULL%1,               Mattia%Busana%NULL%2,               Mattia%Busana%NULL%0,               Sandra%Rossi%NULL%0,               Sandra%Rossi%NULL%0,               Davide%Chiumello%NULL%0,               Davide%Chiumello%NULL%0]"
$ws.Range("E4").Value = "[Khai%Tran%NULL%1,               Karen%Cimon%NULL%1,               Melissa%Severn%NULL%1,               Carmem L.%Pessoa-Silva%NULL%1,               John%Conly%NULL%1,               Malcolm Gracie%Semple%NULL%2,               Malcolm Gracie%Semple%NULL%0]"

$wb.Save()
